$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @('D2', '29.363.93'),
    @('E2', '  +0.01%  '),
    @('D3', '1.864.27'),
    @('E3', '  -0.79%  '),
    @('E4', '  -0.10%  '),
    @('D5', '242.93'),
    @('E5', '  -0.04%  '),
    @('D6', '0.7013'),
    @('E6', '  -2.56%  '),
    @('E7', '  -0.02%  '),
    @('D8', '0.07891'),
    @('E8', '  -1.25%  '),
    @('D9', '0.3122'),
    @('E9', '  -0.60%  '),
    @('D10', '24.38'),
    @('E10', '  -2.04%  '),
    @('D11', '0.07787'),
    @('E11', '  -4.52%  '),
    @('D12', '1.867.61'),
    @('E12', '  -0.64%  '),
    @('D13', '5.142'),
    @('E13', '  -1.63%  '),
    @('D14', '92.34'),
    @('E14', '  -2.56%  '),
    @('D15', '0.6974'),
    @('E15', '  -1.78%  '),
    @('D16', '6.521'),
    @('E16', '  +1.83%  '),
    @('D17', '0.000008587'),
    @('E17', '  +1.52%  '),
    @('D18', '29.386.44'),
    @('E18', '  +0.12%  '),
    @('D19', '248.36'),
    @('E19', '  +0.35%  '),
    @('D20', '2.123.00'),
    @('E20', '  +0.39%  '),
    @('D21', '12.98'),
    @('E21', '  -2.41%  '),
    @('D22', '0.9999'),
    @('E22', '  -0.20%  '),
    @('D23', '7.573'),
    @('E23', '  -2.10%  '),
    @('D24', '1.001'),
    @('E24', '  -0.16%  '),
    @('E25', '  -4.05%  '),
    @('D26', '8.967'),
    @('E26', '  -1.06%  '),
    @('D27', '160.42'),
    @('E27', '  -1.42%  '),
    @('D28', '18.70'),
    @('E28', '  -0.64%  '),
    @('D29', '1.575'),
    @('E29', '  +4.66%  '),
    @('D30', '4.290'),
    @('E30', '  -2.75%  '),
    @('D31', '4.242'),
    @('E31', '  -1.06%  '),
    @('D32', '1.198'),
    @('E32', '  -1.70%  '),
    @('D33', '0.05245'),
    @('E33', '  -1.81%  '),
    @('D34', '1.887'),
    @('E34', '  -2.55%  '),
    @('D35', '0.7563'),
    @('E35', '  +0.01%  '),
    @('D36', '1.181'),
    @('E36', '  +0.25%  '),
    @('D37', '2.712'),
    @('E37', '  +0.35%  '),
    @('D38', '1.272.86'),
    @('E38', '  +0.57%  '),
    @('D39', '0.01864'),
    @('E39', '  -0.93%  '),
    @('D40', '2.749'),
    @('E40', '  -0.41%  '),
    @('D41', '0.8965'),
    @('E41', '  -1.12%  '),
    @('D42', '109.78'),
    @('E42', '  -3.10%  '),
    @('D43', '5.944'),
    @('E43', '  -7.65%  '),
    @('D44', '70.13'),
    @('E44', '  -5.81%  '),
    @('E45', '  -0.11%  '),
    @('D46', '2.020.35'),
    @('E46', '  +0.03%  '),
    @('D47', '0.00000000125'),
    @('E47', '  -3.53%  '),
    @('D48', '9.574'),
    @('E48', '  +0.90%  '),
    @('E49', '  -0.54%  '),
    @('D50', '0.5174'),
    @('E50', '  -0.48%  '),
    @('D51', '0.4282'),
    @('E51', '  -1.39%  '),
)

foreach ($u in $updates) {
    $addr = $u[0]
    $val = $u[1]
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Write-Host "Applied $($updates.Count) cell updates."
